$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 15668
$ws1.Range("F9").Value = 15431
$ws1.Range("F11").Value = 9029
$ws1.Range("F12").Value = 384
$ws1.Range("F15").Value = 91
$ws1.Range("F16").Value = 199
$ws1.Range("F35").Value = 254
$ws1.Range("F36").Value = 325
$ws1.Range("F39").Value = 5554

# Sheet 4: "全部类型" (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 15668
$ws4.Range("F9").Value = 15431
$ws4.Range("F11").Value = 9029
$ws4.Range("F12").Value = 384
$ws4.Range("F15").Value = 91
$ws4.Range("F16").Value = 199
$ws4.Range("F37").Value = 254
$ws4.Range("F38").Value = 325
$ws4.Range("F41").Value = 5554
